$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6941904991429927
$ws.Range("C4").Value = 0.736
$ws.Range("D4").Value = 0.6704154704990131
$ws.Range("E4").Value = 0.681
$ws.Range("F4").Value = 0.5806642707733458
$ws.Range("G4").Value = 0.591
$ws.Range("H4").Value = 0.5912946686249376
$ws.Range("I4").Value = 0.5840000000000001
$ws.Range("J4").Value = 0.6514788551601871
$ws.Range("K4").Value = 0.8380000000000001
$ws.Range("L4").Value = 0.5412345486117429
$ws.Range("M4").Value = 0.5615
$ws.Range("N4").Value = 0.6680543586551942
$ws.Range("O4").Value = 0.6729999999999999
$ws.Range("P4").Value = 0.688965916006648
$ws.Range("Q4").Value = 0.6765000000000001
$ws.Range("R4").Value = 0.6947598556616859
$ws.Range("S4").Value = 0.727
$ws.Range("T4").Value = 0.6793016443676528
$ws.Range("U4").Value = 0.6859999999999999
$ws.Range("V4").Value = 0.5771931644746148
$ws.Range("W4").Value = 0.588
$ws.Range("X4").Value = 0.5866105649488949
$ws.Range("Y4").Value = 0.5805
$ws.Range("Z4").Value = 0.6639549072055098
$ws.Range("AA4").Value = 0.669
$ws.Range("AB4").Value = 0.6868777782898029
$ws.Range("AC4").Value = 0.674
$ws.Range("B5").Value = 0.6739828158237947
$ws.Range("C5").Value = 0.7010000000000001
$ws.Range("D5").Value = 0.6591649210188556
$ws.Range("E5").Value = 0.6645
$ws.Range("F5").Value = 0.6369963308024753
$ws.Range("G5").Value = 0.626
$ws.Range("H5").Value = 0.6772249246268817
$ws.Range("I5").Value = 0.6555
$ws.Range("J5").Value = 0.6098216461575732
$ws.Range("K5").Value = 0.744
$ws.Range("L5").Value = 0.5350906230144894
$ws.Range("M5").Value = 0.5464999999999999
$ws.Range("N5").Value = 0.649625947523073
$ws.Range("P5").Value = 0.6764594154676955
$ws.Range("Q5").Value = 0.6649999999999999
$ws.Range("R5").Value = 0.6718838056414631
$ws.Range("S5").Value = 0.704
$ws.Range("T5").Value = 0.6529666935454165
$ws.Range("U5").Value = 0.6609999999999999
$ws.Range("V5").Value = 0.6379756763366037
$ws.Range("W5").Value = 0.626
$ws.Range("X5").Value = 0.6788742392842588
$ws.Range("Y5").Value = 0.657
$ws.Range("Z5").Value = 0.6480326560393256
$ws.Range("AA5").Value = 0.65
$ws.Range("AB5").Value = 0.6749273978863711
$ws.Range("AC5").Value = 0.6635
$ws.Range("B6").Value = 0.6855517295295155
$ws.Range("C6").Value = 0.7110000000000001
$ws.Range("D6").Value = 0.6758777895739467
$ws.Range("E6").Value = 0.6799999999999999
$ws.Range("F6").Value = 0.6394481958969979
$ws.Range("G6").Value = 0.628
$ws.Range("H6").Value = 0.6721551013536382
$ws.Range("I6").Value = 0.6529999999999999
$ws.Range("J6").Value = 0.6484509885942709
$ws.Range("K6").Value = 0.8400000000000001
$ws.Range("L6").Value = 0.5340110659913415
$ws.Range("M6").Value = 0.5515
$ws.Range("N6").Value = 0.6613982280863564
$ws.Range("O6").Value = 0.665
$ws.Range("P6").Value = 0.6877253822548266
$ws.Range("Q6").Value = 0.674
$ws.Range("R6").Value = 0.6906265504642914
$ws.Range("S6").Value = 0.708
$ws.Range("T6").Value = 0.6905939282052941
$ws.Range("U6").Value = 0.6900000000000001
$ws.Range("V6").Value = 0.643734078671087
$ws.Range("W6").Value = 0.636
$ws.Range("X6").Value = 0.6724398638045888
$ws.Range("Y6").Value = 0.6545000000000001
$ws.Range("Z6").Value = 0.6631759228750086
$ws.Range("AA6").Value = 0.667
$ws.Range("AB6").Value = 0.6899740816968076
$ws.Range("AC6").Value = 0.6765
$ws.Range("B7").Value = 0.6098077193257156
$ws.Range("C7").Value = 0.7020000000000001
$ws.Range("D7").Value = 0.5473492129556523
$ws.Range("E7").Value = 0.5569999999999999
$ws.Range("F7").Value = 0.4996096323036269
$ws.Range("G7").Value = 0.512
$ws.Range("H7").Value = 0.5040707999692915
$ws.Range("I7").Value = 0.5039999999999999
$ws.Range("J7").Value = 0.5236761504601554
$ws.Range("K7").Value = 0.545
$ws.Range("L7").Value = 0.528719582101294
$ws.Range("M7").Value = 0.5290000000000001
$ws.Range("N7").Value = 0.5929986052169279
$ws.Range("O7").Value = 0.6339999999999999
$ws.Range("P7").Value = 0.5975920509060326
$ws.Range("Q7").Value = 0.5845
$ws.Range("R7").Value = 0.618869654927676
$ws.Range("S7").Value = 0.705
$ws.Range("T7").Value = 0.5605612444297328
$ws.Range("U7").Value = 0.5725
$ws.Range("V7").Value = 0.4990582438048378
$ws.Range("W7").Value = 0.511
$ws.Range("X7").Value = 0.5039555347435449
$ws.Range("Y7").Value = 0.5039999999999999
$ws.Range("Z7").Value = 0.5967142004039145
$ws.Range("AA7").Value = 0.635
$ws.Range("AB7").Value = 0.6024647237952518
$ws.Range("AC7").Value = 0.5894999999999999
